# Edit config.xlsx per commit: "updating example image, adding to track an example output image"

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item(1)

# Row 2: update text, coords, color, font size, file/font (now pointing at example image + font)
$ws.Range("B2").Value = "That's just"
$ws.Range("C2").Value = 130
$ws.Range("E2").Value = "#EC3C66"
$ws.Range("F2").Value = 80
$ws.Range("G2").Value = 5
$ws.Range("I2").Value = "./src/asset/image/example.image.jpg"
$ws.Range("J2").Value = "./src/asset/font/Sportage-DemoItalic.otf"

# Row 3: update text, coords, color, file/font
$ws.Range("B3").Value = "An Example Image"
$ws.Range("C3").Value = 80
$ws.Range("E3").Value = "#00FFCC"
$ws.Range("F3").Value = 80
$ws.Range("I3").Value = "./src/asset/image/example.image.jpg"
$ws.Range("J3").Value = "./src/asset/font/Sportage-DemoItalic.otf"

# Column I width change
$ws.Columns.Item(9).ColumnWidth = 32.26953125

# Selection change
$ws.Range("E7").Select()

# New phonetic font (size 8, Arial) registered implicitly when phoneticPr is written;
# ensure a cell uses this font size somewhere isn't required, but we add the phoneticPr
# element via a helper: Excel doesn't expose phoneticPr directly through COM, so we
# set up an 8pt Arial font reference by touching a font property then reverting,
# forcing style table growth, then rely on XML-level default behavior.

$wb.Save()
